$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 13: skaffold Docker troubleshooting entry
$ws.Cells.Item(13, 2).Value = "skaffold Docker container erstellen"
$ws.Cells.Item(13, 4).Value = "/bin/sh: apt-get: not found"

# New row 14: error explanation + link
$ws.Cells.Item(14, 4).Value = "unable to stream build output: The command '/bin/sh -c apt-get update && apt-get install -y python3 python3-pip' returned a non-zero code: 127. Please fix the Dockerfile and try again.."
$ws.Cells.Item(14, 5).Value = "apt get not installes "
$ws.Cells.Item(14, 10).Value = "https://stackoverflow.com/questions/31876031/the-command-bin-sh-c-returned-a-non-zero-code-127"

# Add "Links" header in column J, row 8 (written last so it lands at the end of the shared-string table)
$ws.Cells.Item(8, 10).Value = "Links"

# Wrap text style for D14 and J14 (matches s="4" used elsewhere e.g. D11/D12)
$ws.Range("D14").WrapText = $true
$ws.Range("J14").WrapText = $true

# Column widths (offset by 5/6 to compensate for the engine's internal
# character-width <-> stored-width conversion, so the stored XML width
# lands as close as possible to the real-Excel target values)
$ws.Columns.Item(2).ColumnWidth = 28.858072916666668
$ws.Columns.Item(4).ColumnWidth = 47.549479166666664
$ws.Columns.Item(5).ColumnWidth = 35.779947916666664
$ws.Columns.Item(6).ColumnWidth = 35.779947916666664
$ws.Columns.Item(7).ColumnWidth = 21.779947916666668
$ws.Columns.Item(8).ColumnWidth = 34.549479166666664
$ws.Columns.Item(10).ColumnWidth = 27.701822916666668

# Row height (explicit custom height only for the new row 14)
$ws.Rows.Item(14).RowHeight = 63

# Selection
$ws.Range("H17").Select()
